{"js": "// Update the worksheet date line and every \"a OP b = c\" answer cell in the\n// practice table to the new values from the commit.\n// Each entry is [oldText, newText]; every oldText occurs exactly once in\n// the document (the date paragraph + 100 table cells), so a search +\n// whole-match replace per pair reproduces the diff exactly, run-property\n// (font/size) included, with no extra runs created.\nconst pairs = [\n  [\"2025-03-01 Saturday\", \"2025-03-02 Sunday\"],\n  [\"21+33=54\", \"29+49=78\"],\n  [\"72+6=78\", \"39+27=66\"],\n  [\"84-30=54\", \"36-3=33\"],\n  [\"92+3=95\", \"35+56=91\"],\n  [\"56-8=48\", \"40+47=87\"],\n  [\"92-81=11\", \"69+27=96\"],\n  [\"33+29=62\", \"23+20=43\"],\n  [\"3+32=35\", \"22-4=18\"],\n  [\"12+18=30\", \"29+2=31\"],\n  [\"86-39=47\", \"59+5=64\"],\n  [\"49-26=23\", \"58-10=48\"],\n  [\"76+12=88\", \"93-44=49\"],\n  [\"19+38=57\", \"74+13=87\"],\n  [\"91-86=5\", \"69-17=52\"],\n  [\"61+26=87\", \"41-0=41\"],\n  [\"24+58=82\", \"24-11=13\"],\n  [\"68-61=7\", \"47-19=28\"],\n  [\"64+25=89\", \"99-77=22\"],\n  [\"0+59=59\", \"45+30=75\"],\n  [\"52-31=21\", \"20+33=53\"],\n  [\"93-62=31\", \"45+40=85\"],\n  [\"3+15=18\", \"53-14=39\"],\n  [\"35+30=65\", \"91+4=95\"],\n  [\"38+12=50\", \"64-25=39\"],\n  [\"47+40=87\", \"80+9=89\"],\n  [\"26+35=61\", \"71-68=3\"],\n  [\"44+26=70\", \"26+54=80\"],\n  [\"8+69=77\", \"31+12=43\"],\n  [\"13+21=34\", \"93-21=72\"],\n  [\"92-25=67\", \"33+9=42\"],\n  [\"17+74=91\", \"91-77=14\"],\n  [\"99-18=81\", \"70+2=72\"],\n  [\"38+31=69\", \"72-65=7\"],\n  [\"35-4=31\", \"94+2=96\"],\n  [\"42+15=57\", \"1+35=36\"],\n  [\"90+9=99\", \"63-60=3\"],\n  [\"58-6=52\", \"33+59=92\"],\n  [\"90-66=24\", \"88+6=94\"],\n  [\"57-22=35\", \"11+13=24\"],\n  [\"2+3=5\", \"44-34=10\"],\n  [\"72+23=95\", \"6+86=92\"],\n  [\"0+75=75\", \"73-71=2\"],\n  [\"94-79=15\", \"10+61=71\"],\n  [\"28-23=5\", \"84-0=84\"],\n  [\"38+14=52\", \"92-35=57\"],\n  [\"26+16=42\", \"55-19=36\"],\n  [\"89-62=27\", \"5+63=68\"],\n  [\"81-78=3\", \"34+15=49\"],\n  [\"83-78=5\", \"58+14=72\"],\n  [\"83-72=11\", \"90-89=1\"],\n  [\"13+48=61\", \"80-75=5\"],\n  [\"13+45=58\", \"13-9=4\"],\n  [\"54+44=98\", \"99-35=64\"],\n  [\"73-57=16\", \"61-6=55\"],\n  [\"66+29=95\", \"98-43=55\"],\n  [\"31+11=42\", \"98-73=25\"],\n  [\"34+27=61\", \"21+73=94\"],\n  [\"35-5=30\", \"37-15=22\"],\n  [\"80+18=98\", \"13+24=37\"],\n  [\"77-9=68\", \"69+22=91\"],\n  [\"6-0=6\", \"71+17=88\"],\n  [\"47-6=41\", \"95-32=63\"],\n  [\"95-51=44\", \"80+12=92\"],\n  [\"30+42=72\", \"27+61=88\"],\n  [\"67-6=61\", \"42+44=86\"],\n  [\"4+68=72\", \"56-54=2\"],\n  [\"55-31=24\", \"49-35=14\"],\n  [\"14+85=99\", \"21-5=16\"],\n  [\"7+68=75\", \"27+35=62\"],\n  [\"5-0=5\", \"28+36=64\"],\n  [\"30+7=37\", \"83-29=54\"],\n  [\"41-35=6\", \"37-31=6\"],\n  [\"35+46=81\", \"15+41=56\"],\n  [\"50-44=6\", \"45+22=67\"],\n  [\"53+45=98\", \"33-20=13\"],\n  [\"47+20=67\", \"93-54=39\"],\n  [\"70-27=43\", \"76-64=12\"],\n  [\"29+5=34\", \"60+18=78\"],\n  [\"65+27=92\", \"91-61=30\"],\n  [\"66-54=12\", \"55+22=77\"],\n  [\"54-33=21\", \"86-50=36\"],\n  [\"57+19=76\", \"13-9=4\"],\n  [\"2+59=61\", \"91-20=71\"],\n  [\"19+13=32\", \"49+34=83\"],\n  [\"29+60=89\", \"47-41=6\"],\n  [\"71+5=76\", \"61-58=3\"],\n  [\"36-12=24\", \"35+12=47\"],\n  [\"43-10=33\", \"12+70=82\"],\n  [\"55+30=85\", \"60-47=13\"],\n  [\"79-20=59\", \"53-34=19\"],\n  [\"30-7=23\", \"13+76=89\"],\n  [\"6+89=95\", \"1+2=3\"],\n  [\"40-19=21\", \"32+37=69\"],\n  [\"37+55=92\", \"15-5=10\"],\n  [\"83-5=78\", \"84-58=26\"],\n  [\"18+58=76\", \"11+30=41\"],\n  [\"71+19=90\", \"77-44=33\"],\n  [\"57-25=32\", \"70+22=92\"],\n  [\"82-71=11\", \"24+26=50\"],\n  [\"69-0=69\", \"35+39=74\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  // Replace in place so the run's formatting (font/size) is preserved.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date line and every \"a OP b = c\" answer cell in the\n# practice table to the new values from the commit.\n# Each entry is (oldText, newText); every oldText occurs exactly once in\n# the document (the date paragraph + 100 table cells), so Find/Replace\n# (MatchCase on, one occurrence at a time) reproduces the diff exactly,\n# leaving run formatting (font/size) untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"2025-03-01 Saturday\", \"2025-03-02 Sunday\")\n    ,@(\"21+33=54\", \"29+49=78\")\n    ,@(\"72+6=78\", \"39+27=66\")\n    ,@(\"84-30=54\", \"36-3=33\")\n    ,@(\"92+3=95\", \"35+56=91\")\n    ,@(\"56-8=48\", \"40+47=87\")\n    ,@(\"92-81=11\", \"69+27=96\")\n    ,@(\"33+29=62\", \"23+20=43\")\n    ,@(\"3+32=35\", \"22-4=18\")\n    ,@(\"12+18=30\", \"29+2=31\")\n    ,@(\"86-39=47\", \"59+5=64\")\n    ,@(\"49-26=23\", \"58-10=48\")\n    ,@(\"76+12=88\", \"93-44=49\")\n    ,@(\"19+38=57\", \"74+13=87\")\n    ,@(\"91-86=5\", \"69-17=52\")\n    ,@(\"61+26=87\", \"41-0=41\")\n    ,@(\"24+58=82\", \"24-11=13\")\n    ,@(\"68-61=7\", \"47-19=28\")\n    ,@(\"64+25=89\", \"99-77=22\")\n    ,@(\"0+59=59\", \"45+30=75\")\n    ,@(\"52-31=21\", \"20+33=53\")\n    ,@(\"93-62=31\", \"45+40=85\")\n    ,@(\"3+15=18\", \"53-14=39\")\n    ,@(\"35+30=65\", \"91+4=95\")\n    ,@(\"38+12=50\", \"64-25=39\")\n    ,@(\"47+40=87\", \"80+9=89\")\n    ,@(\"26+35=61\", \"71-68=3\")\n    ,@(\"44+26=70\", \"26+54=80\")\n    ,@(\"8+69=77\", \"31+12=43\")\n    ,@(\"13+21=34\", \"93-21=72\")\n    ,@(\"92-25=67\", \"33+9=42\")\n    ,@(\"17+74=91\", \"91-77=14\")\n    ,@(\"99-18=81\", \"70+2=72\")\n    ,@(\"38+31=69\", \"72-65=7\")\n    ,@(\"35-4=31\", \"94+2=96\")\n    ,@(\"42+15=57\", \"1+35=36\")\n    ,@(\"90+9=99\", \"63-60=3\")\n    ,@(\"58-6=52\", \"33+59=92\")\n    ,@(\"90-66=24\", \"88+6=94\")\n    ,@(\"57-22=35\", \"11+13=24\")\n    ,@(\"2+3=5\", \"44-34=10\")\n    ,@(\"72+23=95\", \"6+86=92\")\n    ,@(\"0+75=75\", \"73-71=2\")\n    ,@(\"94-79=15\", \"10+61=71\")\n    ,@(\"28-23=5\", \"84-0=84\")\n    ,@(\"38+14=52\", \"92-35=57\")\n    ,@(\"26+16=42\", \"55-19=36\")\n    ,@(\"89-62=27\", \"5+63=68\")\n    ,@(\"81-78=3\", \"34+15=49\")\n    ,@(\"83-78=5\", \"58+14=72\")\n    ,@(\"83-72=11\", \"90-89=1\")\n    ,@(\"13+48=61\", \"80-75=5\")\n    ,@(\"13+45=58\", \"13-9=4\")\n    ,@(\"54+44=98\", \"99-35=64\")\n    ,@(\"73-57=16\", \"61-6=55\")\n    ,@(\"66+29=95\", \"98-43=55\")\n    ,@(\"31+11=42\", \"98-73=25\")\n    ,@(\"34+27=61\", \"21+73=94\")\n    ,@(\"35-5=30\", \"37-15=22\")\n    ,@(\"80+18=98\", \"13+24=37\")\n    ,@(\"77-9=68\", \"69+22=91\")\n    ,@(\"6-0=6\", \"71+17=88\")\n    ,@(\"47-6=41\", \"95-32=63\")\n    ,@(\"95-51=44\", \"80+12=92\")\n    ,@(\"30+42=72\", \"27+61=88\")\n    ,@(\"67-6=61\", \"42+44=86\")\n    ,@(\"4+68=72\", \"56-54=2\")\n    ,@(\"55-31=24\", \"49-35=14\")\n    ,@(\"14+85=99\", \"21-5=16\")\n    ,@(\"7+68=75\", \"27+35=62\")\n    ,@(\"5-0=5\", \"28+36=64\")\n    ,@(\"30+7=37\", \"83-29=54\")\n    ,@(\"41-35=6\", \"37-31=6\")\n    ,@(\"35+46=81\", \"15+41=56\")\n    ,@(\"50-44=6\", \"45+22=67\")\n    ,@(\"53+45=98\", \"33-20=13\")\n    ,@(\"47+20=67\", \"93-54=39\")\n    ,@(\"70-27=43\", \"76-64=12\")\n    ,@(\"29+5=34\", \"60+18=78\")\n    ,@(\"65+27=92\", \"91-61=30\")\n    ,@(\"66-54=12\", \"55+22=77\")\n    ,@(\"54-33=21\", \"86-50=36\")\n    ,@(\"57+19=76\", \"13-9=4\")\n    ,@(\"2+59=61\", \"91-20=71\")\n    ,@(\"19+13=32\", \"49+34=83\")\n    ,@(\"29+60=89\", \"47-41=6\")\n    ,@(\"71+5=76\", \"61-58=3\")\n    ,@(\"36-12=24\", \"35+12=47\")\n    ,@(\"43-10=33\", \"12+70=82\")\n    ,@(\"55+30=85\", \"60-47=13\")\n    ,@(\"79-20=59\", \"53-34=19\")\n    ,@(\"30-7=23\", \"13+76=89\")\n    ,@(\"6+89=95\", \"1+2=3\")\n    ,@(\"40-19=21\", \"32+37=69\")\n    ,@(\"37+55=92\", \"15-5=10\")\n    ,@(\"83-5=78\", \"84-58=26\")\n    ,@(\"18+58=76\", \"11+30=41\")\n    ,@(\"71+19=90\", \"77-44=33\")\n    ,@(\"57-25=32\", \"70+22=92\")\n    ,@(\"82-71=11\", \"24+26=50\")\n    ,@(\"69-0=69\", \"35+39=74\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    # wdFindContinue=1, wdReplaceOne=1\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 1) | Out-Null\n}\n"}
